$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3708.9443
$ws.Range("I19").Value = 3674.6667
$ws.Range("J19").Value = 3726.0833
$ws.Range("K19").Value = 3674.6667
$ws.Range("L19").Value = 3726.0833
$ws.Range("M19").Value = -3499.6667
$ws.Range("N19").Value = -4076.0833
$ws.Range("H40").Value = 15991.643
$ws.Range("I40").Value = 9999.833000000001
$ws.Range("K40").Value = 9999.833000000001
$ws.Range("M40").Value = -9824.833000000001
$ws.Range("H96").Value = 1555
$ws.Range("I96").Value = 1000.0909
$ws.Range("K96").Value = 3000.2727
$ws.Range("M96").Value = -1627.2727
$ws.Range("H98").Value = 1385.5294
$ws.Range("I98").Value = 1200.875
$ws.Range("K98").Value = 1200.875
$ws.Range("M98").Value = 297.125
$ws.Range("H122").Value = 1385.5294
$ws.Range("I122").Value = 1200.875
$ws.Range("K122").Value = 3602.625
$ws.Range("M122").Value = -1152.625
$ws.Range("H137").Value = 2311.598
$ws.Range("I137").Value = 1943.1111
$ws.Range("J137").Value = 2835.2368
$ws.Range("K137").Value = 5829.3333
$ws.Range("L137").Value = 8505.7104
$ws.Range("M137").Value = -3279.3333
$ws.Range("N137").Value = -13605.7104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 189973
$ws.Range("J69").Value = 189973
$ws.Range("L69").Value = 189973
$ws.Range("N69").Value = -191471
$ws.Range("H72").Value = 189973
$ws.Range("J72").Value = 189973
$ws.Range("L72").Value = 569919
$ws.Range("N72").Value = -577407
$ws.Range("H74").Value = 1476.1154
$ws.Range("I74").Value = 1489.1765
$ws.Range("J74").Value = 1451.4445
$ws.Range("K74").Value = 1489.1765
$ws.Range("L74").Value = 1451.4445
$ws.Range("M74").Value = -615.1765
$ws.Range("N74").Value = -3199.4445
$ws.Range("H77").Value = 1476.1154
$ws.Range("I77").Value = 1489.1765
$ws.Range("J77").Value = 1451.4445
$ws.Range("K77").Value = 7445.8825
$ws.Range("L77").Value = 7257.2225
$ws.Range("M77").Value = -3077.8825
$ws.Range("N77").Value = -15993.2225
$ws.Range("H88").Value = 1472.25
$ws.Range("I88").Value = 1236
$ws.Range("J88").Value = 1614
$ws.Range("K88").Value = 1236
$ws.Range("L88").Value = 1614
$ws.Range("M88").Value = -830
$ws.Range("N88").Value = -2426
$ws.Range("H91").Value = 1472.25
$ws.Range("I91").Value = 1236
$ws.Range("J91").Value = 1614
$ws.Range("K91").Value = 1236
$ws.Range("L91").Value = 1614
$ws.Range("M91").Value = 168
$ws.Range("N91").Value = -4422
$ws.Range("H102").Value = 3153.6
$ws.Range("I102").Value = 2605.862
$ws.Range("K102").Value = 2605.862
$ws.Range("M102").Value = -983.8620000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 99999
$ws.Range("J42").Value = 99999
$ws.Range("L42").Value = 99999
$ws.Range("N42").Value = -100655
$ws.Range("H43").Value = 99999
$ws.Range("J43").Value = 99999
$ws.Range("L43").Value = 99999
$ws.Range("N43").Value = -100361
$ws.Range("H48").Value = 110341.5
$ws.Range("J48").Value = 110341.5
$ws.Range("L48").Value = 110341.5
$ws.Range("N48").Value = -111171.5
$ws.Range("H70").Value = 169822
$ws.Range("J70").Value = 169822
$ws.Range("L70").Value = 169822
$ws.Range("N70").Value = -170408
$ws.Range("H73").Value = 169822
$ws.Range("J73").Value = 169822
$ws.Range("L73").Value = 169822
$ws.Range("N73").Value = -171850
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H134").Value = 4454.2554
$ws.Range("I134").Value = 3086.1025
$ws.Range("K134").Value = 9258.307499999999
$ws.Range("M134").Value = -6723.307499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2459.9185
$ws.Range("I31").Value = 2324.1724
$ws.Range("J31").Value = 2656.75
$ws.Range("K31").Value = 2324.1724
$ws.Range("L31").Value = 2656.75
$ws.Range("M31").Value = -2029.1724
$ws.Range("N31").Value = -3246.75
$ws.Range("H34").Value = 2459.9185
$ws.Range("I34").Value = 2324.1724
$ws.Range("J34").Value = 2656.75
$ws.Range("K34").Value = 2324.1724
$ws.Range("L34").Value = 2656.75
$ws.Range("M34").Value = -2122.1724
$ws.Range("N34").Value = -3060.75
$ws.Range("H122").Value = 12263.462
$ws.Range("I122").Value = 13724.412
$ws.Range("K122").Value = 41173.236
$ws.Range("M122").Value = -38723.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 14929
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 14929
$ws.Range("K101").Value = 0
$ws.Range("M101").Value = 44787
$ws.Range("N101").Value = -49655
$ws.Range("H113").Value = 849.6429000000001
$ws.Range("J113").Value = 731.25
$ws.Range("L113").Value = 2193.75
$ws.Range("N113").Value = -6533.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4168.9287
$ws.Range("I70").Value = 3069.75
$ws.Range("J70").Value = 5634.5
$ws.Range("K70").Value = 3069.75
$ws.Range("L70").Value = 5634.5
$ws.Range("M70").Value = -2799.75
$ws.Range("N70").Value = -6174.5
$ws.Range("H73").Value = 4168.9287
$ws.Range("I73").Value = 3069.75
$ws.Range("J73").Value = 5634.5
$ws.Range("K73").Value = 3069.75
$ws.Range("L73").Value = 5634.5
$ws.Range("M73").Value = -2133.75
$ws.Range("N73").Value = -7506.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3284.7878
$ws.Range("I40").Value = 2690.4517
$ws.Range("J40").Value = 12497
$ws.Range("K40").Value = 2690.4517
$ws.Range("L40").Value = 12497
$ws.Range("M40").Value = -2554.4517
$ws.Range("N40").Value = -12769
$ws.Range("H46").Value = 2527.9092
$ws.Range("I46").Value = 1051.8334
$ws.Range("J46").Value = 3081.4375
$ws.Range("K46").Value = 1051.8334
$ws.Range("L46").Value = 3081.4375
$ws.Range("M46").Value = -863.8334
$ws.Range("N46").Value = -3457.4375
$ws.Range("H136").Value = 7833692.5
$ws.Range("I136").Value = 13849384
$ws.Range("K136").Value = 41548152
$ws.Range("M136").Value = -41545602
$ws.Range("H141").Value = 242999.25
$ws.Range("J141").Value = 242999.25
$ws.Range("L141").Value = 242999.25
$ws.Range("N141").Value = -253359.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4467.25
$ws.Range("I126").Value = 4064.7307
$ws.Range("J126").Value = 5513.8
$ws.Range("K126").Value = 12194.1921
$ws.Range("L126").Value = 16541.4
$ws.Range("M126").Value = -9724.1921
$ws.Range("N126").Value = -21481.4
$ws.Range("H136").Value = 23016.105
$ws.Range("I136").Value = 23037.75
$ws.Range("J136").Value = 22900.666
$ws.Range("K136").Value = 69113.25
$ws.Range("L136").Value = 68701.99800000001
$ws.Range("M136").Value = -66563.25
$ws.Range("N136").Value = -73801.99800000001
